$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 66 (bus 64) entirely - the sheet now only has buses 0..63
$ws.Rows(66).Delete()

# Build the updated B2:E65 block (vm_pu, va_degree, p_mw, q_mvar) with the new solved values
$data = New-Object 'object[,]' 64,4
$data[0,0] = [double]"1.000000000000023"
$data[0,1] = [double]"0"
$data[0,2] = [double]"12.9518182153464"
$data[0,3] = [double]"-13.91692037604958"
$data[1,0] = [double]"1.000599999999979"
$data[1,1] = [double]"0.1494997929456079"
$data[1,2] = [double]"-6.804032892659961e-10"
$data[1,3] = [double]"-81.82204788025294"
$data[2,0] = [double]"1.000599999999979"
$data[2,1] = [double]"0.1494997929456079"
$data[2,2] = [double]"-6.804032892659961e-10"
$data[2,3] = [double]"-81.82204788025294"
$data[3,0] = [double]"1.000000000000016"
$data[3,1] = [double]"0.2189475394097234"
$data[3,2] = [double]"-8.609224444455776e-13"
$data[3,3] = [double]"43.59119939151923"
$data[4,0] = [double]"1.000450150037185"
$data[4,1] = [double]"0.1714224146058856"
$data[4,2] = [double]"-0.04899999999648562"
$data[4,3] = [double]"-38.62972680574185"
$data[5,0] = [double]"0.9999695787768388"
$data[5,1] = [double]"0.227499900948209"
$data[5,2] = [double]"-1.960029361036675e-12"
$data[5,3] = [double]"20.12892978089425"
$data[6,0] = [double]"0.9999272134206989"
$data[6,1] = [double]"0.2345810896443885"
$data[6,2] = [double]"-6.55160647955455e-12"
$data[6,3] = [double]"14.33169189917605"
$data[7,0] = [double]"1.000172277403482"
$data[7,1] = [double]"0.2084712104063883"
$data[7,2] = [double]"-0.04130000000444532"
$data[7,3] = [double]"-18.01463148657798"
$data[8,0] = [double]"0.9999821177386103"
$data[8,1] = [double]"0.2306649767693044"
$data[8,2] = [double]"2.526329145879913e-11"
$data[8,3] = [double]"24.12684108810313"
$data[9,0] = [double]"1.000389001114892"
$data[9,1] = [double]"0.1860384030537106"
$data[9,2] = [double]"-0.01644999999885582"
$data[9,3] = [double]"-7.794703266869609"
$data[10,0] = [double]"1.000607987694504"
$data[10,1] = [double]"0.1623238610025166"
$data[10,2] = [double]"-0.01960000001322637"
$data[10,3] = [double]"-32.32307561087495"
$data[11,0] = [double]"1.000049327941521"
$data[11,1] = [double]"0.2252999286310448"
$data[11,2] = [double]"8.478218127550008e-13"
$data[11,3] = [double]"21.99854282538907"
$data[12,0] = [double]"1.000023752017013"
$data[12,1] = [double]"0.2286327818246822"
$data[12,2] = [double]"-0.03639999999487682"
$data[12,3] = [double]"-4.110341096483543"
$data[13,0] = [double]"1.000015092464526"
$data[13,1] = [double]"0.2055176110382032"
$data[13,2] = [double]"-0.02940000000018941"
$data[13,3] = [double]"-1.885533507370745"
$data[14,0] = [double]"1.000018831733736"
$data[14,1] = [double]"0.2067939188901947"
$data[14,2] = [double]"-0.03009999999644071"
$data[14,3] = [double]"-5.152037274573039"
$data[15,0] = [double]"0.9999774616846271"
$data[15,1] = [double]"0.2491530662307657"
$data[15,2] = [double]"-0.04900000000297691"
$data[15,3] = [double]"2.162383814530362"
$data[16,0] = [double]"0.9999777186814519"
$data[16,1] = [double]"0.2619853934001002"
$data[16,2] = [double]"-0.04829999999852769"
$data[16,3] = [double]"1.595506407077778"
$data[17,0] = [double]"1.000001632011149"
$data[17,1] = [double]"0.254295076238436"
$data[17,2] = [double]"-0.04900000000042661"
$data[17,3] = [double]"-0.6636235824793493"
$data[18,0] = [double]"0.9999749788342462"
$data[18,1] = [double]"0.2394384057055946"
$data[18,2] = [double]"-0.04199999999836908"
$data[18,3] = [double]"0.7469843808462976"
$data[19,0] = [double]"1.0000534779554"
$data[19,1] = [double]"0.2278015259952622"
$data[19,2] = [double]"-0.04199999999962995"
$data[19,3] = [double]"0.209927636398003"
$data[20,0] = [double]"0.9999201244174717"
$data[20,1] = [double]"0.2907369126554896"
$data[20,2] = [double]"-0.04200000000125986"
$data[20,3] = [double]"5.31388484284335"
$data[21,0] = [double]"1.000599999999979"
$data[21,1] = [double]"0.1494997929456079"
$data[21,2] = [double]"-6.804032892659961e-10"
$data[21,3] = [double]"-81.82204788025294"
$data[22,0] = [double]"0.9999999999999706"
$data[22,1] = [double]"0.8478386835868799"
$data[22,2] = [double]"-2.591629320002065e-10"
$data[22,3] = [double]"0.2559975470168256"
$data[23,0] = [double]"0.9997699007235318"
$data[23,1] = [double]"1.297447993508358"
$data[23,2] = [double]"-7.000000000001045"
$data[23,3] = [double]"23.99790724370353"
$data[24,0] = [double]"0.9997790088075567"
$data[24,1] = [double]"1.090715067037111"
$data[24,2] = [double]"7.862183126761124e-13"
$data[24,3] = [double]"1.443378940324558"
$data[25,0] = [double]"1.000599999999979"
$data[25,1] = [double]"0.1494997929456079"
$data[25,2] = [double]"-6.804032892659961e-10"
$data[25,3] = [double]"-81.82204788025294"
$data[26,0] = [double]"1.000000000000007"
$data[26,1] = [double]"0.304761245419105"
$data[26,2] = [double]"-3.288229064035697e-11"
$data[26,3] = [double]"4.171529408926449"
$data[27,0] = [double]"0.9997005556993197"
$data[27,1] = [double]"0.4310159557377004"
$data[27,2] = [double]"-0.400000000022091"
$data[27,3] = [double]"3.788409001081372"
$data[28,0] = [double]"0.9996237209868081"
$data[28,1] = [double]"0.5292198549459598"
$data[28,2] = [double]"-0.400000000022027"
$data[28,3] = [double]"8.045616502456669"
$data[29,0] = [double]"1.000298092425488"
$data[29,1] = [double]"0.5241280801138205"
$data[29,2] = [double]"-0.4000000000031588"
$data[29,3] = [double]"-9.425870143577132"
$data[30,0] = [double]"1.000259107248875"
$data[30,1] = [double]"0.5832245894232876"
$data[30,2] = [double]"-0.4000000000033745"
$data[30,3] = [double]"2.159019654020798"
$data[31,0] = [double]"1.000312369854602"
$data[31,1] = [double]"0.6152957510578961"
$data[31,2] = [double]"-0.4000000000041276"
$data[31,3] = [double]"4.041312441530417"
$data[32,0] = [double]"1.000602215656174"
$data[32,1] = [double]"0.6175373933199906"
$data[32,2] = [double]"-0.399999999996291"
$data[32,3] = [double]"-18.39828403208596"
$data[33,0] = [double]"0.9995479719766764"
$data[33,1] = [double]"0.7466538489808388"
$data[33,2] = [double]"-0.400000000000728"
$data[33,3] = [double]"20.24564163668481"
$data[34,0] = [double]"1.000024771832486"
$data[34,1] = [double]"0.7101849857409915"
$data[34,2] = [double]"-0.3999999999980175"
$data[34,3] = [double]"-13.68278629377882"
$data[35,0] = [double]"0.9996328653931489"
$data[35,1] = [double]"0.6424312904389202"
$data[35,2] = [double]"-0.4000000000074668"
$data[35,3] = [double]"13.265799659922"
$data[36,0] = [double]"0.9995727952059271"
$data[36,1] = [double]"0.6864842938454327"
$data[36,2] = [double]"1.648098324480429e-12"
$data[36,3] = [double]"9.413262923251832"
$data[37,0] = [double]"0.9996900346867543"
$data[37,1] = [double]"0.7169767090012458"
$data[37,2] = [double]"-0.4000000000003133"
$data[37,3] = [double]"0.9972776178983512"
$data[38,0] = [double]"1.000240558966568"
$data[38,1] = [double]"0.4526999614208178"
$data[38,2] = [double]"-1.182776099284411e-12"
$data[38,3] = [double]"-8.158931199493836"
$data[39,0] = [double]"1.000786197009713"
$data[39,1] = [double]"0.36380527360739"
$data[39,2] = [double]"-0.7000000000005783"
$data[39,3] = [double]"-9.845792622249455"
$data[40,0] = [double]"1.000408788243711"
$data[40,1] = [double]"0.6064733958461273"
$data[40,2] = [double]"-0.4000000000056491"
$data[40,3] = [double]"-10.19371130817999"
$data[41,0] = [double]"1.000475545840634"
$data[41,1] = [double]"0.589959181221587"
$data[41,2] = [double]"-6.65155430734643e-13"
$data[41,3] = [double]"-5.233707496080896"
$data[42,0] = [double]"0.9995848131452133"
$data[42,1] = [double]"0.9557670942943588"
$data[42,2] = [double]"-0.400000000001079"
$data[42,3] = [double]"12.12598671278865"
$data[43,0] = [double]"1.00011238767836"
$data[43,1] = [double]"0.7146047212546496"
$data[43,2] = [double]"1.278491201794907e-13"
$data[43,3] = [double]"4.353433053434671"
$data[44,0] = [double]"0.9995529982161124"
$data[44,1] = [double]"0.8333862666466098"
$data[44,2] = [double]"-0.4999999999999307"
$data[44,3] = [double]"2.941877203918968"
$data[45,0] = [double]"0.9994972100534627"
$data[45,1] = [double]"0.9429736650571167"
$data[45,2] = [double]"-0.3000000000003436"
$data[45,3] = [double]"8.072087075574862"
$data[46,0] = [double]"1.000599999999979"
$data[46,1] = [double]"0.1494997929456079"
$data[46,2] = [double]"-6.804032892659961e-10"
$data[46,3] = [double]"-81.82204788025294"
$data[47,0] = [double]"1.000000000000008"
$data[47,1] = [double]"0.2178213278370268"
$data[47,2] = [double]"-2.574412905076429e-12"
$data[47,3] = [double]"53.95129712733974"
$data[48,0] = [double]"1.000700000000009"
$data[48,1] = [double]"0.1413860833182438"
$data[48,2] = [double]"1.792158088598228e-11"
$data[48,3] = [double]"-69.52168862738318"
$data[49,0] = [double]"1.000399999999992"
$data[49,1] = [double]"0.1754493700758464"
$data[49,2] = [double]"-0.03499999998742297"
$data[49,3] = [double]"16.68759424836658"
$data[50,0] = [double]"1.000500000000009"
$data[50,1] = [double]"0.1647078246772885"
$data[50,2] = [double]"-0.03000000000192792"
$data[50,3] = [double]"-16.57869745945102"
$data[51,0] = [double]"1.000199999999993"
$data[51,1] = [double]"0.1984030011424019"
$data[51,2] = [double]"-0.03999999998665889"
$data[51,3] = [double]"10.46062432047728"
$data[52,0] = [double]"1.000149999999993"
$data[52,1] = [double]"0.2039794320610423"
$data[52,2] = [double]"-8.974524248750448e-12"
$data[52,3] = [double]"2.077724659856653"
$data[53,0] = [double]"1.000029999999998"
$data[53,1] = [double]"0.217396558005239"
$data[53,2] = [double]"-0.04500000001099712"
$data[53,3] = [double]"28.97220568135782"
$data[54,0] = [double]"1.000056000000001"
$data[54,1] = [double]"0.2156145239100115"
$data[54,2] = [double]"-0.0649999999878438"
$data[54,3] = [double]"-0.1643522222293929"
$data[55,0] = [double]"1.000075000000016"
$data[55,1] = [double]"0.2141841576576106"
$data[55,2] = [double]"-0.01499999999945605"
$data[55,3] = [double]"-0.03314925640285285"
$data[56,0] = [double]"1.000090000000007"
$data[56,1] = [double]"0.2130167193517191"
$data[56,2] = [double]"-0.04999999999982768"
$data[56,3] = [double]"-3.954892315516527"
$data[57,0] = [double]"1.000010000000006"
$data[57,1] = [double]"0.2222029201407074"
$data[57,2] = [double]"6.204897706751922e-12"
$data[57,3] = [double]"17.81928797989706"
$data[58,0] = [double]"1.00002699999999"
$data[58,1] = [double]"0.2205300004574725"
$data[58,2] = [double]"-0.01225000001446733"
$data[58,3] = [double]"11.48585686172193"
$data[59,0] = [double]"1.000320000000002"
$data[59,1] = [double]"0.1880135855794825"
$data[59,2] = [double]"-0.01749999999825178"
$data[59,3] = [double]"-22.60606326131886"
$data[60,0] = [double]"1.000104495842711"
$data[60,1] = [double]"0.2121251527024888"
$data[60,2] = [double]"-0.01539999999738114"
$data[60,3] = [double]"10.47389947409463"
$data[61,0] = [double]"1.000294595884538"
$data[61,1] = [double]"0.190513673819049"
$data[61,2] = [double]"-0.01225000000863635"
$data[61,3] = [double]"-17.69148763178461"
$data[62,0] = [double]"1.000214030884966"
$data[62,1] = [double]"0.1995091704515161"
$data[62,2] = [double]"-0.01434999999899684"
$data[62,3] = [double]"21.41763169736491"
$data[63,0] = [double]"1.000574859636533"
$data[63,1] = [double]"0.1592356265305944"
$data[63,2] = [double]"-1.801428450853848e-11"
$data[63,3] = [double]"-17.51413970337486"

$ws.Range("B2:E65").Value = $data
